$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 6 (ano = 2025) with refreshed metrics
$ws.Range("C6").Value = 348
$ws.Range("D6").Value = 281
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 62.58351893095768
$ws.Range("G6").Value = 19.25287356321839
$ws.Range("H6").Value = 80.74712643678161
